# Applies the "Pushing newer QDesc Jutsu" edit:
#  - Row 2 (Age) stats are replaced with new values
#  - Row 3 label changes from "Salary" to "sleep(hrs)/weeknights" and its
#    stats are replaced with new values
#  - Row 4 (the old "Satisfaction" row) is removed entirely
#  - Sheet dimension shrinks from A1:K4 to A1:K3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 ("Age") statistics
$ws.Range("B2").Value = 671
$ws.Range("C2").Value = 49.8
$ws.Range("D2").Value = 19.36
$ws.Range("E2").Value = 49
$ws.Range("F2").Value = 17
$ws.Range("G2").Value = 25.2
$ws.Range("H2").Value = 18
$ws.Range("I2").Value = 84
$ws.Range("J2").Value = 7.38
$ws.Range("K2").Value = 0.01

# Update row 3 label and statistics
$ws.Range("A3").Value = "sleep(hrs)/weeknights"
$ws.Range("B3").Value = 671
$ws.Range("C3").Value = 7.63
$ws.Range("D3").Value = 1.39
$ws.Range("E3").Value = 8.34
$ws.Range("F3").Value = 0.98
$ws.Range("G3").Value = 1.45
$ws.Range("H3").Value = 4.52
$ws.Range("I3").Value = 9.880000000000001
$ws.Range("J3").Value = 35.7
$ws.Range("K3").Value = 0.01

# Remove row 4 (old "Satisfaction" row) entirely
$ws.Rows.Item(4).Delete()
